$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Partial Molar Volume Analysis")

$ws.Range("F2").Value = 0.21681059625562235
$ws.Range("N2").Value = 52.46213447313791
$ws.Range("V2").Value = 0.21681059625562235
$ws.Range("Z2").Value = 52.46213447313791
$ws.Range("D3").Value = 2.0282703317762607
$ws.Range("F3").Value = 0.21681059625562235
$ws.Range("N3").Value = 52.47485911759988
$ws.Range("T3").Value = 0.4052995299870725
$ws.Range("V3").Value = 0.21681059625562235
$ws.Range("Z3").Value = 52.69507199388379
$ws.Range("D4").Value = 4.224784285543178
$ws.Range("F4").Value = 0.21681059625562235
$ws.Range("N4").Value = 52.48895451574734
$ws.Range("T4").Value = 0.810599059974145
$ws.Range("V4").Value = 0.21681059625562235
$ws.Range("Z4").Value = 52.92799553235364
$ws.Range("D5").Value = 5.981995448556714
$ws.Range("F5").Value = 0.21681059625562235
$ws.Range("N5").Value = 52.500466783677986
$ws.Range("T5").Value = 1.2158985899612178
$ws.Range("V5").Value = 0.21681059625562235
$ws.Range("Z5").Value = 53.16090508966752
$ws.Range("D6").Value = 7.944837705114385
$ws.Range("F6").Value = 0.21681059625562235
$ws.Range("N6").Value = 52.513574140654804
$ws.Range("T6").Value = 1.62119811994829
$ws.Range("V6").Value = 0.21681059625562235
$ws.Range("Z6").Value = 53.393800666942234
$ws.Range("D7").Value = 9.888986225895316
$ws.Range("F7").Value = 0.21681059625562235
$ws.Range("N7").Value = 52.52681447785474
$ws.Range("T7").Value = 2.0264976499353624
$ws.Range("V7").Value = 0.21681059625562235
$ws.Range("Z7").Value = 53.62668226529891
$ws.Range("D8").Value = 11.753686369625104
$ws.Range("F8").Value = 0.21681059625562235
$ws.Range("N8").Value = 52.539754737584474
$ws.Range("T8").Value = 2.4317971799224356
$ws.Range("V8").Value = 0.21681059625562235
$ws.Range("Z8").Value = 53.859549885855444
$ws.Range("T9").Value = 2.837096709909508
$ws.Range("V9").Value = 0.21681059625562235
$ws.Range("Z9").Value = 54.09240352972972
$ws.Range("T10").Value = 3.24239623989658
$ws.Range("V10").Value = 0.21681059625562235
$ws.Range("Z10").Value = 54.32524319803964
$ws.Range("T11").Value = 3.6476957698836525
$ws.Range("V11").Value = 0.21681059625562235
$ws.Range("Z11").Value = 54.55806889190634
$ws.Range("T12").Value = 4.052995299870725
$ws.Range("V12").Value = 0.21681059625562235
$ws.Range("Z12").Value = 54.790880612445534
$ws.Range("T13").Value = 4.458294829857798
$ws.Range("V13").Value = 0.21681059625562235
$ws.Range("Z13").Value = 55.02367836077728
$ws.Range("T14").Value = 4.863594359844871
$ws.Range("V14").Value = 0.21681059625562235
$ws.Range("Z14").Value = 55.2564621380184
$ws.Range("T15").Value = 5.2688938898319435
$ws.Range("V15").Value = 0.21681059625562235
$ws.Range("Z15").Value = 55.4892319452857
$ws.Range("T16").Value = 5.674193419819016
$ws.Range("V16").Value = 0.21681059625562235
$ws.Range("Z16").Value = 55.7219877837003
$ws.Range("T17").Value = 6.079492949806087
$ws.Range("V17").Value = 0.21681059625562235
$ws.Range("Z17").Value = 55.9547296543758
$ws.Range("T18").Value = 6.48479247979316
$ws.Range("V18").Value = 0.21681059625562235
$ws.Range("Z18").Value = 56.18745755843222
$ws.Range("T19").Value = 6.890092009780233
$ws.Range("V19").Value = 0.21681059625562235
$ws.Range("Z19").Value = 56.4201714969864
$ws.Range("T20").Value = 7.295391539767305
$ws.Range("V20").Value = 0.21681059625562235
$ws.Range("Z20").Value = 56.6528714711562
$ws.Range("T21").Value = 7.700691069754377
$ws.Range("V21").Value = 0.21681059625562235
$ws.Range("Z21").Value = 56.88555748205955
$ws.Range("T22").Value = 8.10599059974145
$ws.Range("V22").Value = 0.21681059625562235
$ws.Range("Z22").Value = 57.118229530812165
$ws.Range("T23").Value = 8.511290129728524
$ws.Range("V23").Value = 0.21681059625562235
$ws.Range("Z23").Value = 57.350887618530855
$ws.Range("T24").Value = 8.916589659715596
$ws.Range("V24").Value = 0.21681059625562235
$ws.Range("Z24").Value = 57.58353174633459
$ws.Range("T25").Value = 9.321889189702668
$ws.Range("V25").Value = 0.21681059625562235
$ws.Range("Z25").Value = 57.816161915338036
$ws.Range("T26").Value = 9.727188719689742
$ws.Range("V26").Value = 0.21681059625562235
$ws.Range("Z26").Value = 58.04877812665908
$ws.Range("T27").Value = 10.132488249676813
$ws.Range("V27").Value = 0.21681059625562235
$ws.Range("Z27").Value = 58.28138038141345
$ws.Range("T28").Value = 10.537787779663887
$ws.Range("V28").Value = 0.21681059625562235
$ws.Range("Z28").Value = 58.51396868071689
$ws.Range("T29").Value = 10.943087309650958
$ws.Range("V29").Value = 0.21681059625562235
$ws.Range("Z29").Value = 58.74654302568944
$ws.Range("T30").Value = 11.348386839638032
$ws.Range("V30").Value = 0.21681059625562235
$ws.Range("Z30").Value = 58.97910341744361
$ws.Range("T31").Value = 11.753686369625104
$ws.Range("V31").Value = 0.21681059625562235
$ws.Range("Z31").Value = 59.21164985709727
